# Refresh the cryptos price table (GitHub Actions scheduled update).
#
# Column D ("Price") and E ("Volume(1h)") hold plain-text values even when
# they look numeric (e.g. "22.10", "0.9989"), matching the original file's
# inline-string cells. A leading apostrophe forces Excel to keep the
# assigned value as text instead of silently re-parsing it as a number
# (which would drop trailing zeros, e.g. "22.10" -> 22.1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.475.54'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '1.648.38'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("D4").Value = '''0.9989'
$ws.Range("E4").Value = '  -0.59%  '
$ws.Range("D5").Value = '''0.9998'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = '''300.44'
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("D8").Value = '''50.38'
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("D9").Value = '''0.3506'
$ws.Range("E9").Value = '  -2.63%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '''1.223'
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.08081'
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").Value = '''0.9989'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = '''22.10'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = '''6.330'
$ws.Range("E14").Value = '  -1.52%  '
$ws.Range("D15").Value = '''7.268'
$ws.Range("E15").Value = '  -2.07%  '
$ws.Range("D16").Value = '''0.00001218'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").Value = '1.647.26'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '''94.96'
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("D19").Value = '''0.06969'
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("D20").Value = '''6.629'
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("D21").Value = '''17.46'
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = '''0.9993'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").Value = '''12.46'
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("D24").Value = '23.476.63'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = '''2.435'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").Value = '''2.995'
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("D27").Value = '''21.10'
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("D28").Value = '''150.66'
$ws.Range("E28").Value = '  -1.51%  '
$ws.Range("D29").Value = '''5.186'
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").Value = '''131.75'
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("D31").Value = '1.832.45'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").Value = '''6.890'
$ws.Range("E32").Value = '  -2.10%  '
$ws.Range("D33").Value = '''2.133'
$ws.Range("E33").Value = '  -5.40%  '
$ws.Range("D34").Value = '''11.19'
$ws.Range("E34").Value = '  -8.15%  '
$ws.Range("D35").Value = '''0.9909'
$ws.Range("E35").Value = '  -5.89%  '
$ws.Range("D36").Value = '''0.02702'
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("D37").Value = '''0.08785'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E38").Value = '  -2.70%  '
$ws.Range("D39").Value = '''5.916'
$ws.Range("E39").Value = '  -2.27%  '
$ws.Range("D40").Value = '''0.06802'
$ws.Range("E40").Value = '  -2.32%  '
$ws.Range("D41").Value = '''12.82'
$ws.Range("E41").Value = '  -1.98%  '
$ws.Range("D42").Value = '''0.6864'
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("D43").Value = '''1.293'
$ws.Range("E43").Value = '  -2.94%  '
$ws.Range("D44").Value = '''15.45'
$ws.Range("E44").Value = '  -3.02%  '
$ws.Range("D45").Value = '''0.9986'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").Value = '''0.6380'
$ws.Range("E46").Value = '  -1.46%  '
$ws.Range("D47").Value = '''2.249'
$ws.Range("E47").Value = '  -1.58%  '
$ws.Range("D48").Value = '''3.919'
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("D49").Value = '''0.07681'
$ws.Range("E49").Value = '  -2.39%  '
$ws.Range("D50").Value = '''126.98'
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("D51").Value = '''1.229'
$ws.Range("E51").Value = '  +2.61%  '
